$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the current row 2 (pushes all existing
# trade rows down by 3 -> old row 2 becomes row 5, etc.)
$ws.Range("A2:A4").EntireRow.Insert()

# New row 2: BTCUSDT BUY
$ws.Cells.Item(2,1).Value = 45931.03042824074
$ws.Cells.Item(2,2).Value = "BTCUSDT"
$ws.Cells.Item(2,3).Value = "BUY"
$ws.Cells.Item(2,4).Value = 113954.4
$ws.Cells.Item(2,5).Value = 0.00396
$ws.Cells.Item(2,6).Value = 451.259424
$ws.Cells.Item(2,7).Value = 0.00003834
$ws.Cells.Item(2,8).Value = "BNB"

# New row 3: TRXUSDT SELL (fee coin USDC - new shared string)
$ws.Cells.Item(3,1).Value = 45931.02980324074
$ws.Cells.Item(3,2).Value = "TRXUSDT"
$ws.Cells.Item(3,3).Value = "SELL"
$ws.Cells.Item(3,4).Value = 0.3329
$ws.Cells.Item(3,5).Value = 905.7
$ws.Cells.Item(3,6).Value = 301.50753
$ws.Cells.Item(3,7).Value = 0.28643215
$ws.Cells.Item(3,8).Value = "USDC"

# New row 4: TRXUSDT SELL (fee coin USDC)
$ws.Cells.Item(4,1).Value = 45931.028958333336
$ws.Cells.Item(4,2).Value = "TRXUSDT"
$ws.Cells.Item(4,3).Value = "SELL"
$ws.Cells.Item(4,4).Value = 0.3329
$ws.Cells.Item(4,5).Value = 452
$ws.Cells.Item(4,6).Value = 150.4708
$ws.Cells.Item(4,7).Value = 0.14294726
$ws.Cells.Item(4,8).Value = "USDC"

# Move the selection to match the saved cursor position after editing
$ws.Range("F4").Select()
